# Auto-generated script updating the cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force text storage so number-like strings (e.g. '11.60', '0.0000143')
    # are not auto-converted to numeric cells, then restore the default
    # (unstyled) cell style so no stray formatting is introduced.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '70.567.14'
$ws.Range('E2').Value = '  +2.03%  '

$ws.Range('D3').Value = '3.801.81'
$ws.Range('E3').Value = '  +0.63%  '

$ws.Range('E4').Value = '  +0.04%  '

Set-TextValue 'D5' '667.23'
$ws.Range('E5').Value = '  +6.72%  '

Set-TextValue 'D6' '168.99'
$ws.Range('E6').Value = '  +1.66%  '

$ws.Range('D7').Value = '3.798.86'
$ws.Range('E7').Value = '  +0.59%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('E9').Value = '  +0.96%  '

$ws.Range('E10').Value = '  +0.15%  '

$ws.Range('E11').Value = '  +1.49%  '

Set-TextValue 'D12' '7.03'
$ws.Range('E12').Value = '  +4.74%  '

$ws.Range('E13').Value = '  -0.88%  '

Set-TextValue 'D14' '35.72'
$ws.Range('E14').Value = '  +0.24%  '

$ws.Range('D15').Value = '4.442.98'

$ws.Range('D16').Value = '3.805.19'
$ws.Range('E16').Value = '  +0.40%  '

$ws.Range('D17').Value = '70.507.41'
$ws.Range('E17').Value = '  +1.94%  '

Set-TextValue 'D18' '17.69'
$ws.Range('E18').Value = '  +0.19%  '

Set-TextValue 'D19' '7.18'
$ws.Range('E19').Value = '  +1.02%  '

Set-TextValue 'D20' '11.60'
$ws.Range('E20').Value = '  +20.41%  '

$ws.Range('E21').Value = '  +0.54%  '

Set-TextValue 'D22' '474.19'
$ws.Range('E22').Value = '  +1.24%  '

Set-TextValue 'D23' '0.713'
$ws.Range('E23').Value = '  +0.99%  '

Set-TextValue 'D24' '82.95'
$ws.Range('E24').Value = '  -0.24%  '

Set-TextValue 'D25' '0.0000143'
$ws.Range('E25').Value = '  -2.98%  '

Set-TextValue 'D26' '12.20'
$ws.Range('E26').Value = '  +1.29%  '

Set-TextValue 'D27' '10.33'
$ws.Range('E27').Value = '  +3.16%  '

Set-TextValue 'D28' '2.12'
$ws.Range('E28').Value = '  -1.84%  '

$ws.Range('E29').Value = '  +0.01%  '

$ws.Range('D30').Value = '3.953.74'
$ws.Range('E30').Value = '  +0.70%  '

Set-TextValue 'D31' '2.85'
$ws.Range('E31').Value = '  +6.68%  '

$ws.Range('E32').Value = '  +2.80%  '

Set-TextValue 'D33' '7.42'
$ws.Range('E33').Value = '  +2.43%  '

Set-TextValue 'D34' '29.57'
$ws.Range('E34').Value = '  +2.78%  '

$ws.Range('E35').Value = '  +6.47%  '

Set-TextValue 'D36' '9.11'
$ws.Range('E36').Value = '  +1.29%  '

$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('D38').Value = '3.759.19'
$ws.Range('E38').Value = '  +0.83%  '

Set-TextValue 'D39' '0.103'
$ws.Range('E39').Value = '  +0.94%  '

Set-TextValue 'D40' '3.41'
$ws.Range('E40').Value = '  -0.07%  '

Set-TextValue 'D41' '5.95'
$ws.Range('E41').Value = '  +2.53%  '

Set-TextValue 'D42' '0.965'
$ws.Range('E42').Value = '  -0.42%  '

$ws.Range('E43').Value = '  +0.10%  '

$ws.Range('E44').Value = '  +10.07%  '

Set-TextValue 'D46' '45.64'
$ws.Range('E46').Value = '  +6.11%  '

Set-TextValue 'D47' '158.76'
$ws.Range('E47').Value = '  +4.08%  '

Set-TextValue 'D48' '47.99'
$ws.Range('E48').Value = '  +2.82%  '

$ws.Range('E49').Value = '  +5.04%  '

Set-TextValue 'D50' '0.299'
$ws.Range('E50').Value = '  +0.47%  '

$ws.Range('E51').Value = '  +1.32%  '

